# Testing permissions, - test example 9
#
# Renames several "permissions" scenario-step strings in the
# "ATDD Scenarios" sheet (Given-When-Then (Description), column H) of
# Table29, and moves the active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios")

# --- H column (Given-When-Then (Description)) text updates -----------------
# "Non-restrictive base permissions"  -> "Full base starting permissions" / "Full base permissions"
# "Non-restrictive base permissions extended with Lookup Value permissions"
#      -> "Full base starting permissions extended with Lookup Value permissions" / "Full base permissions extended with Lookup Value"
# "Unrestricted permissions" -> "Unrestricted starting permissions"

$ws.Range("H171").Value = "Full base starting permissions"
$ws.Range("H175").Value = "Full base starting permissions extended with Lookup Value permissions"
$ws.Range("H179").Value = "Unrestricted starting permissions"
$ws.Range("H181").Value = "Full base permissions"
$ws.Range("H185").Value = "Unrestricted starting permissions"
$ws.Range("H187").Value = "Full base permissions extended with Lookup Value"
$ws.Range("H191").Value = "Unrestricted starting permissions"
$ws.Range("H193").Value = "Full base permissions"
$ws.Range("H197").Value = "Unrestricted starting permissions"
$ws.Range("H199").Value = "Full base permissions extended with Lookup Value"
$ws.Range("H203").Value = "Unrestricted starting permissions"
$ws.Range("H205").Value = "Full base permissions"
$ws.Range("H209").Value = "Unrestricted starting permissions"
$ws.Range("H211").Value = "Full base permissions extended with Lookup Value"
$ws.Range("H215").Value = "Unrestricted starting permissions"
$ws.Range("H217").Value = "Full base permissions"
$ws.Range("H221").Value = "Full base permissions extended with Lookup Value"
$ws.Range("H226").Value = "Full base starting permissions"
$ws.Range("H230").Value = "Full base starting permissions extended with Lookup Value permissions"
$ws.Range("H234").Value = "Full base starting permissions"
$ws.Range("H238").Value = "Full base starting permissions extended with Lookup Value permissions"

# Re-fit the rows whose wrapped-text line count didn't change, so the
# edit doesn't leave a stray explicit row height behind.
$ws.Rows("171:171").AutoFit()
$ws.Rows("175:175").AutoFit()
$ws.Rows("179:179").AutoFit()
$ws.Rows("181:181").AutoFit()
$ws.Rows("185:185").AutoFit()
$ws.Rows("187:187").AutoFit()
$ws.Rows("191:191").AutoFit()
$ws.Rows("193:193").AutoFit()
$ws.Rows("197:197").AutoFit()
$ws.Rows("199:199").AutoFit()
$ws.Rows("203:203").AutoFit()
$ws.Rows("205:205").AutoFit()
$ws.Rows("209:209").AutoFit()
$ws.Rows("211:211").AutoFit()
$ws.Rows("215:215").AutoFit()
$ws.Rows("217:217").AutoFit()
$ws.Rows("221:221").AutoFit()
$ws.Rows("226:226").AutoFit()
$ws.Rows("234:234").AutoFit()

# Rows 230 and 238 now wrap onto a second line, so Excel grows them to
# the two-line height.
$ws.Rows("230:230").RowHeight = 30
$ws.Rows("238:238").RowHeight = 30

# --- Selection change --------------------------------------------------
$ws.Range("B127").Select()
